$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New document row: insert a blank row at row 26 (shifts the existing last
# row - "29/7/2025(Onsite)" / "Car Tracking Project" / the cleaner-steps
# note - down to row 27), leaving room to add the new doc entry.
$ws.Rows.Item(26).Insert()

# Reflect the view-state Excel recorded after the edit (scroll position / zoom).
$win = $excel.ActiveWindow
$win.ScrollRow = 18
$win.ScrollColumn = 1
$win.Zoom = 126
$ws.Range("C26").Select()

Write-Output "Inserted blank row at 26; row 26 content shifted to row 27."
